$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 0
$ws.Range("F6").Value = 700848
$ws.Range("F7").Value = 1965
$ws.Range("F11").Value = 1403
$ws.Range("F12").Value = 1233
$ws.Range("F13").Value = 2918
$ws.Range("F14").Value = 1736
$ws.Range("F15").Value = 1136
$ws.Range("F17").Value = 50
$ws.Range("F18").Value = 11
$ws.Range("F19").Value = 16
$ws.Range("F20").Value = 580
$ws.Range("I22").Value = "//i2.hdslb.com/bfs/openplatform/202407/BeA1vFig1720509195913.jpeg"
$ws.Range("F23").Value = 1712
$ws.Range("I23").Value = "//i2.hdslb.com/bfs/openplatform/202407/BeA1vFig1720509195913.jpeg"
$ws.Range("F24").Value = 0
$ws.Range("F25").Value = 159
$ws.Range("F26").Value = 0
$ws.Range("F27").Value = 1651
$ws.Range("F28").Value = 631
$ws.Range("F29").Value = 0
$ws.Range("F30").Value = 38
$ws.Range("F32").Value = 1197
$ws.Range("F33").Value = 114
$ws.Range("F35").Value = 221
$ws.Range("F36").Value = 349
$ws.Range("F37").Value = 0
$ws.Range("F40").Value = 1096
$ws.Range("F41").Value = 0
$ws.Range("F42").Value = 1085
$ws.Range("F43").Value = 42
$ws.Range("F44").Value = 910
$ws.Range("F45").Value = 211
$ws.Range("F46").Value = 699
$ws.Range("F48").Value = 58
$ws.Range("F50").Value = 0
$ws.Range("F51").Value = 48

# --- Sheet: 演出 ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 4
$ws.Range("F5").Value = 0
$ws.Range("F8").Value = 148222
$ws.Range("F10").Value = 64
$ws.Range("F11").Value = 30
$ws.Range("F14").Value = 243
$ws.Range("F15").Value = 369
$ws.Range("F17").Value = 448
$ws.Range("F18").Value = 1
$ws.Range("F19").Value = 237
$ws.Range("F21").Value = 102
$ws.Range("F22").Value = 95
$ws.Range("F23").Value = 0
$ws.Range("F25").Value = 1
$ws.Range("F29").Value = 0
$ws.Range("F31").Value = 0
$ws.Range("F32").Value = 126
$ws.Range("F33").Value = 126
$ws.Range("F34").Value = 26
$ws.Range("F36").Value = 257
$ws.Range("F37").Value = 116
$ws.Range("F39").Value = 197
$ws.Range("F41").Value = 15
$ws.Range("F43").Value = 8
$ws.Range("F44").Value = 0

# --- Sheet: 本地生活 ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 0
$ws.Range("F4").Value = 3192
$ws.Range("F5").Value = 303
$ws.Range("F8").Value = 0
$ws.Range("F9").Value = 697
$ws.Range("F11").Value = 2340

# --- Sheet: 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 0
$ws.Range("F4").Value = 172
$ws.Range("F5").Value = 252
$ws.Range("F6").Value = 266
$ws.Range("F7").Value = 0
$ws.Range("F11").Value = 148222
$ws.Range("F12").Value = 1677
$ws.Range("F13").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("F17").Value = 2918
$ws.Range("F19").Value = 0
$ws.Range("F20").Value = 0
$ws.Range("F21").Value = 1690
$ws.Range("F22").Value = 11
$ws.Range("F26").Value = 1166
$ws.Range("F27").Value = 1712
$ws.Range("I27").Value = "//i2.hdslb.com/bfs/openplatform/202407/BeA1vFig1720509195913.jpeg"
$ws.Range("F28").Value = 1712
$ws.Range("I28").Value = "//i2.hdslb.com/bfs/openplatform/202407/BeA1vFig1720509195913.jpeg"
$ws.Range("F29").Value = 1177
$ws.Range("F30").Value = 159
$ws.Range("F32").Value = 631
$ws.Range("F34").Value = 38
$ws.Range("F35").Value = 1197
$ws.Range("F36").Value = 0
$ws.Range("F37").Value = 465
$ws.Range("F40").Value = 13
$ws.Range("F41").Value = 0
$ws.Range("F42").Value = 26
$ws.Range("F43").Value = 349
$ws.Range("F45").Value = 235
$ws.Range("F46").Value = 1097
$ws.Range("F47").Value = 1085
$ws.Range("F48").Value = 910
$ws.Range("F49").Value = 211
$ws.Range("F50").Value = 699
$ws.Range("F51").Value = 58
$ws.Range("F52").Value = 733
